# "with neo4j fix and INS"
# - Update the Neo4j "StatQuery" (CasesTab) query text in cell B2 of the
#   "startup" sheet to also return demo.weight (WITH ... AS age, demo.weight as weight)
# - Move the active selection from B2 to C2
# - Grow row 2's height from 300 to 315 (to fit the longer query text)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$cell = $ws.Range("B2")
$oldText = $cell.Value()

$oldFragment = "WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age"
$newFragment = "WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight"

if ($oldText.Contains($oldFragment)) {
    $newText = $oldText.Replace($oldFragment, $newFragment)
    $cell.Value = $newText
}

# Grow the row to accommodate the extra text.
$ws.Rows.Item(2).RowHeight = 315

# Move the selection to C2.
$ws.Range("C2").Select()
